$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B21: convert stored text "3" into a real number 3 ---
$ws.Range("B21").Value = 3

# --- Add new row 22 ---
$ws.Range("A22").Value = "Ruilin"

# B22 must stay a *text* cell containing "4" (not a number), even though
# it looks numeric. A direct Value assignment of a numeric-looking string
# gets auto-coerced to a number by this engine, so build it as a text
# formula result in a scratch cell, then paste just the value into B22.
$ws.Range("Z1").Formula = '=T("4")'
$ws.Range("Z1").Copy()
$ws.Range("B22").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("C22").Value = "propose a new method"
$ws.Range("D22").Value = "DIS"
$ws.Range("E22").Value = "MET"
$ws.Range("F22").Value = "55e6f9d2-bdcc-4319-8467-87a8dbd0172d"
$ws.Range("G22").Value = "Byt3oJ-0W_annotated.xlsx"
$ws.Range("H22").Value = "The authors propose a new method that approximates the discrete max-weight matching by a continuous Sinkhorn operator, which looks like an analog of softmax operator on matrices."
